$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell value updates (modified or newly added cells)
$ws.Range("K3").Value = 823
$ws.Range("K5").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("H23").Value = 0
$ws.Range("J24").Value = 2206.333333333335
$ws.Range("K26").Value = 931.333333333333
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 1299
$ws.Range("J29").Value = 790.666666666667
$ws.Range("K29").Value = 1331.333333333333
$ws.Range("H31").Value = 1
$ws.Range("F32").Value = 'incongruent'
$ws.Range("G32").Value = 6526
$ws.Range("I32").Value = 1
$ws.Range("J32").Value = 0
$ws.Range("H33").Value = 1
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("H44").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K46").Value = 1357.333333333333
$ws.Range("J48").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("K54").Value = 649
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("H71").Value = 0
$ws.Range("H72").Value = 0
$ws.Range("H75").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("I84").Value = 1
$ws.Range("H87").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("J94").Value = 615.666666666667
$ws.Range("G101").Value = 8042.666666666667
$ws.Range("H101").Value = 1
$ws.Range("J101").Value = 0
$ws.Range("J102").Value = 590.6666666666679
$ws.Range("J105").Value = 0
$ws.Range("J107").Value = 2048
$ws.Range("I114").Value = 2
$ws.Range("I118").Value = 0
$ws.Range("H123").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 1557.333333333334
$ws.Range("J125").Value = 0
$ws.Range("J127").Value = 480.333333333333
$ws.Range("K129").Value = 0
$ws.Range("K137").Value = 865.6666666666661
$ws.Range("H138").Value = 0
$ws.Range("J143").Value = 906.3333333333339
$ws.Range("J145").Value = 622.9999999999991
$ws.Range("K146").Value = 840.666666666667
$ws.Range("K147").Value = 749
$ws.Range("J151").Value = 0
$ws.Range("J152").Value = 0
$ws.Range("I154").Value = 1
$ws.Range("J155").Value = 424
$ws.Range("F156").Value = 'congruent'
$ws.Range("G156").Value = 7851
$ws.Range("H156").Value = 1
$ws.Range("J164").Value = 1206.333333333332
$ws.Range("J165").Value = 1024
$ws.Range("I175").Value = 0
$ws.Range("I176").Value = 2
$ws.Range("H177").Value = 2
$ws.Range("I181").Value = 0
$ws.Range("K195").Value = 0
$ws.Range("H196").Value = 1
$ws.Range("J196").Value = 932.333333333333
$ws.Range("K199").Value = 0
$ws.Range("F200").Value = 'congruent'
$ws.Range("G200").Value = 9609.333333333334
$ws.Range("I200").Value = 0
$ws.Range("I201").Value = 0
$ws.Range("F213").Value = 'incongruent'
$ws.Range("G213").Value = 8867.666666666668
$ws.Range("H213").Value = 0
$ws.Range("G217").Value = 8826
$ws.Range("H217").Value = 1
$ws.Range("K218").Value = 1239.666666666667
$ws.Range("J224").Value = 0
$ws.Range("J225").Value = 157.3333333333339
$ws.Range("J226").Value = 0
$ws.Range("H230").Value = 0
$ws.Range("I231").Value = 0
$ws.Range("G234").Value = 7226
$ws.Range("K236").Value = 0
$ws.Range("J237").Value = 806.333333333333
$ws.Range("J243").Value = 0
$ws.Range("H244").Value = 0
$ws.Range("H245").Value = 2
$ws.Range("G247").Value = 7834.333333333334
$ws.Range("I247").Value = 1
$ws.Range("I249").Value = 0
$ws.Range("H261").Value = 1

# Cells removed entirely (ClearContents removes the <c> element on save)
$ws.Range("F23").ClearContents()
$ws.Range("G23").ClearContents()
$ws.Range("F35").ClearContents()
$ws.Range("G35").ClearContents()
$ws.Range("F44").ClearContents()
$ws.Range("G44").ClearContents()
$ws.Range("F45").ClearContents()
$ws.Range("G45").ClearContents()
$ws.Range("H61").ClearContents()
$ws.Range("I61").ClearContents()
$ws.Range("J62").ClearContents()
$ws.Range("K62").ClearContents()
$ws.Range("F66").ClearContents()
$ws.Range("G66").ClearContents()
$ws.Range("F71").ClearContents()
$ws.Range("G71").ClearContents()
$ws.Range("F72").ClearContents()
$ws.Range("G72").ClearContents()
$ws.Range("F75").ClearContents()
$ws.Range("G75").ClearContents()
$ws.Range("F78").ClearContents()
$ws.Range("G78").ClearContents()
$ws.Range("F87").ClearContents()
$ws.Range("G87").ClearContents()
$ws.Range("F118").ClearContents()
$ws.Range("G118").ClearContents()
$ws.Range("F138").ClearContents()
$ws.Range("G138").ClearContents()
$ws.Range("F181").ClearContents()
$ws.Range("G181").ClearContents()
$ws.Range("F201").ClearContents()
$ws.Range("G201").ClearContents()
$ws.Range("F230").ClearContents()
$ws.Range("G230").ClearContents()
$ws.Range("F231").ClearContents()
$ws.Range("G231").ClearContents()
$ws.Range("F244").ClearContents()
$ws.Range("G244").ClearContents()
$ws.Range("F249").ClearContents()
$ws.Range("G249").ClearContents()
